# "Segundo commit do TI"
#
# The deck originally has a single title slide (slide1.xml, layout
# "Slide de Título" / ctrTitle). This edit appends two new slides, both
# built from the "Título e Conteúdo" (Title and Content) custom layout,
# which is the 2nd layout of the slide master (matches the <p:ph type="title"/>
# + <p:ph idx="1"/> placeholder pair seen in the target slide2.xml/slide3.xml).

$p = $ppt.ActivePresentation

# --- New slide 2: "Trabalho de TI" -----------------------------------
$s2 = $p.Slides.Add(2, 2)

$title2 = $s2.Shapes.Item(1).TextFrame.TextRange
$title2.Text = "Trabalho de TI"
$title2.LanguageID = "pt-BR"

# --- New slide 3: "Slide 3" -------------------------------------------
$s3 = $p.Slides.Add(3, 2)

$title3 = $s3.Shapes.Item(1).TextFrame.TextRange
$title3.Text = "Slide 3"
$title3.LanguageID = "pt-BR"
